# "Actualizacion desde MV -datos-"
# Appends the next batch of daily REPO rows (05-08-2021 .. 03-09-2021) to
# Sheet1, continuing directly after the existing last row (217).
#
# Column A holds the date formatted as dd-mm-yyyy, stored as literal text
# (shared string) exactly like the rest of the column - NOT an Excel date
# serial. Typing the text straight into Range.Value lets the host
# auto-recognise some of these (day <= 12) as dates and convert them to
# serial numbers with a new date number-format, which would also bloat the
# style table. To avoid that, the string is first written as a formula that
# evaluates to text ( ="dd-mm-yyyy" ), then converted in place to a plain
# value via copy / paste-values - this yields a plain shared-string text
# cell with no cell-level style, matching the original rows and leaving
# styles.xml untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("05-08-2021", 6480, 13),
    @("06-08-2021", 6470, 13),
    @("07-08-2021", 6470, 13),
    @("08-08-2021", 6470, 13),
    @("09-08-2021", 6406, 13),
    @("10-08-2021", 6382, 13),
    @("11-08-2021", 6442, 13),
    @("12-08-2021", 6501, 13),
    @("13-08-2021", 6503, 13),
    @("14-08-2021", 6503, 13),
    @("15-08-2021", 6503, 13),
    @("16-08-2021", 6491, 13),
    @("17-08-2021", 6408, 13),
    @("18-08-2021", 6364, 13),
    @("19-08-2021", 6380, 13),
    @("20-08-2021", 6361, 13),
    @("21-08-2021", 6361, 13),
    @("22-08-2021", 6361, 13),
    @("23-08-2021", 6390, 13),
    @("24-08-2021", 6416, 13),
    @("25-08-2021", 6425, 13),
    @("26-08-2021", 6418, 13),
    @("27-08-2021", 6401, 13),
    @("28-08-2021", 6401, 13),
    @("29-08-2021", 6401, 13),
    @("30-08-2021", 6407, 13),
    @("31-08-2021", 6443, 13),
    @("01-09-2021", 6483, 13),
    @("02-09-2021", 6551, 13),
    @("03-09-2021", 6540, 13)
)

$startRow = 218
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $serie = $data[$i][0]
    $monto = $data[$i][1]
    $stock = $data[$i][2]

    $aCell = $ws.Range("A$row")
    $aCell.Formula = '="' + $serie + '"'
    $aCell.Copy()
    $aCell.PasteSpecial(-4163)

    $ws.Range("B$row").Value = $monto
    $ws.Range("C$row").Value = $stock
}
